$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C35 (was "NA" inline string, now empty) while keeping the cell itself present
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = ""
$ws.Range("C35").ClearFormats()

# Add row 36
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "2025-03-31"
$ws.Range("A36").ClearFormats()
$ws.Range("B36").Value = "eaux souterraines"
$ws.Range("C36").Value = 423
$ws.Range("D36").Value = 1

# Add row 37
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "2025-03-31"
$ws.Range("A37").ClearFormats()
$ws.Range("B37").Value = "ruissellement"
$ws.Range("C37").Value = 424
$ws.Range("D37").Value = 4
